# Apply the "Holden scheme" update to the HKL-selection grid on Sheet1.
#
# Net effect (per the target diff):
#   - The 4 rows that used to be labelled "HexGrid-90degTilt{2.5,5,10,15}degRes"
#     (rows 16-19, B column) are renamed to "Holden{2.5,5,10,15}".
#   - The original "HexGrid-90degTilt*degRes" rows are re-appended at the
#     bottom of the table as new rows 20-23 (same style/values as rows 16-19
#     had), extending the used range from A1:W19 to A1:W23.
#   - The column headers in row 2 (C2:M2), which list the 11 "[h, k, l]"
#     triples, are reshuffled into a new order (the trailing "*Pair*"/"5A4F"/
#     "MaxUnique" headers in N2:W2 are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate rows 16-19 (style + values) down into new rows 20-23 *before*
#    renaming row 16-19's labels, so the copies keep the original
#    "HexGrid-90degTilt*degRes" text.
$ws.Range("A16:W19").Copy($ws.Range("A20:W23"))

# 2) Fix up the index column (A) on the newly appended rows to continue the
#    existing 0-based numbering (18, 19, 20, 21).
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21

# 3) Rename the original rows 16-19 to the new "Holden" scheme names.
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# 4) Re-order the "[h, k, l]" column headers in row 2 (C2:M2). The remaining
#    headers (N2:W2 - the pair-grouping names) are unchanged.
$ws.Range("C2").Value = "[5, 1, 1]"
$ws.Range("D2").Value = "[4, 2, 2]"
$ws.Range("E2").Value = "[3, 1, 1]"
$ws.Range("F2").Value = "[3, 3, 1]"
$ws.Range("G2").Value = "[2, 2, 2]"
$ws.Range("H2").Value = "[1, 1, 1]"
$ws.Range("I2").Value = "[3, 3, 3]"
$ws.Range("J2").Value = "[2, 2, 0]"
$ws.Range("K2").Value = "[2, 0, 0]"
$ws.Range("L2").Value = "[4, 0, 0]"
$ws.Range("M2").Value = "[4, 2, 0]"
